# Update "想去人数" (interest count) figures in the F column of the
# "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows keyed by row number on that sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 5556
$ws1.Range("F8").Value  = 911
$ws1.Range("F9").Value  = 146
$ws1.Range("F10").Value = 2488
$ws1.Range("F12").Value = 107
$ws1.Range("F13").Value = 3
$ws1.Range("F14").Value = 72
$ws1.Range("F15").Value = 7
$ws1.Range("F16").Value = 2326
$ws1.Range("F17").Value = 278

# Sheet "全部类型" (same events, but shifted down because it contains
# extra rows from the other category sheets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 5556
$ws4.Range("F10").Value = 911
$ws4.Range("F11").Value = 146
$ws4.Range("F12").Value = 2488
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 3
$ws4.Range("F17").Value = 72
$ws4.Range("F18").Value = 7
$ws4.Range("F19").Value = 2326
$ws4.Range("F20").Value = 278
